# Re-sync row data for rows 2,3,4,5,6,7,11,12 (source Artportalen export re-ordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y (Startdatum) and AA (Slutdatum) are stored as plain text (e.g. "2021-07-20")
# in the source data, not as real dates. Force the target cells to Text format first so
# assigning a date-shaped string does not get auto-converted into a date serial by Excel.
foreach ($addr in @("Y2","AA2","Y4","AA4","Y6","AA6","Y7","AA7","Y11","AA11","Y12","AA12")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 <- source row 11
$ws.Range("A2").Value = 94996219
$ws.Range("B2").Value = 78570
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("P2").Value = "Härjedalen, Hjd"
$ws.Range("Q2").Value = 445370.7671139772
$ws.Range("R2").Value = 6928604.672176878
$ws.Range("S2").Value = 25
$ws.Range("Y2").Value = "2021-07-20"
$ws.Range("AA2").Value = "2021-07-20"
$ws.Range("AC2").Value = "På jättesälg ca 45 cm diameter."
$ws.Range("AW2").Value = "Jens Johannesson"
$ws.Range("AX2").Value = "Jens Johannesson"
$ws.Range("AY2").Value = ""

# Row 3 <- source row 7
$ws.Range("A3").Value = 89596109
$ws.Range("B3").Value = 78570
$ws.Range("E3").Value = 2081
$ws.Range("F3").Value = "Skrovellav"
$ws.Range("G3").Value = "Lobaria scrobiculata"
$ws.Range("H3").Value = "(Scop.) DC."
$ws.Range("Q3").Value = 445260.1075701897
$ws.Range("R3").Value = 6928606.027293501

# Row 4 <- source row 12
$ws.Range("A4").Value = 94996005
$ws.Range("B4").Value = 78570
$ws.Range("E4").Value = 2081
$ws.Range("F4").Value = "Skrovellav"
$ws.Range("G4").Value = "Lobaria scrobiculata"
$ws.Range("H4").Value = "(Scop.) DC."
$ws.Range("P4").Value = "Härjedalen, Hjd"
$ws.Range("Q4").Value = 445261.8150698114
$ws.Range("R4").Value = 6928597.212872105
$ws.Range("S4").Value = 25
$ws.Range("Y4").Value = "2021-07-20"
$ws.Range("AA4").Value = "2021-07-20"
$ws.Range("AC4").Value = "Mkt gammal sälg, fin skog på åsar."
$ws.Range("AW4").Value = "Jens Johannesson"
$ws.Range("AX4").Value = "Jens Johannesson"
$ws.Range("AY4").Value = ""

# Row 5 <- source row 3
$ws.Range("A5").Value = 89596129
$ws.Range("B5").Value = 76909
$ws.Range("E5").Value = 6437
$ws.Range("F5").Value = "Blanksvart spiklav"
$ws.Range("G5").Value = "Calicium denigratum"
$ws.Range("H5").Value = "(Vain.) Tibell"
$ws.Range("Q5").Value = 445032.0268228107
$ws.Range("R5").Value = 6928535.7939387

# Row 6 <- source row 5
$ws.Range("A6").Value = 89596128
$ws.Range("B6").Value = 77506
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("P6").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q6").Value = 445020.0425569176
$ws.Range("R6").Value = 6928539.228960108
$ws.Range("S6").Value = 10
$ws.Range("Y6").Value = "2020-09-25"
$ws.Range("AA6").Value = "2020-09-25"
$ws.Range("AC6").Value = ""
$ws.Range("AW6").Value = "Erland Lindblad"
$ws.Range("AX6").Value = "Jan Henriksson"
$ws.Range("AY6").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 7 <- source row 6
$ws.Range("A7").Value = 94995564
$ws.Range("B7").Value = 95525
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 221941
$ws.Range("F7").Value = "Plattlummer"
$ws.Range("G7").Value = "Lycopodium complanatum"
$ws.Range("H7").Value = "L."
$ws.Range("P7").Value = "Härjedalen, Hjd"
$ws.Range("Q7").Value = 445086.0189850244
$ws.Range("R7").Value = 6928496.057011075
$ws.Range("S7").Value = 25
$ws.Range("Y7").Value = "2021-07-20"
$ws.Range("AA7").Value = "2021-07-20"
$ws.Range("AC7").Value = "Tallskog, ås."
$ws.Range("AW7").Value = "Jens Johannesson"
$ws.Range("AX7").Value = "Jens Johannesson"
$ws.Range("AY7").Value = ""

# Row 11 <- source row 2
$ws.Range("A11").Value = 89596126
$ws.Range("B11").Value = 56411
$ws.Range("E11").Value = 100049
$ws.Range("F11").Value = "Spillkråka"
$ws.Range("G11").Value = "Dryocopus martius"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("P11").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q11").Value = 444929.0050177791
$ws.Range("R11").Value = 6928327.074997591
$ws.Range("S11").Value = 10
$ws.Range("Y11").Value = "2020-09-25"
$ws.Range("AA11").Value = "2020-09-25"
$ws.Range("AC11").Value = "Födosökshack"
$ws.Range("AW11").Value = "Erland Lindblad"
$ws.Range("AX11").Value = "Jan Henriksson"
$ws.Range("AY11").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 12 <- source row 4
$ws.Range("A12").Value = 89596127
$ws.Range("B12").Value = 73693
$ws.Range("E12").Value = 6440
$ws.Range("F12").Value = "Vitgrynig nållav"
$ws.Range("G12").Value = "Chaenotheca subroscida"
$ws.Range("H12").Value = "(Eitner) Zahlbr."
$ws.Range("P12").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q12").Value = 444967.7734563763
$ws.Range("R12").Value = 6928430.952647353
$ws.Range("S12").Value = 10
$ws.Range("Y12").Value = "2020-09-25"
$ws.Range("AA12").Value = "2020-09-25"
$ws.Range("AC12").Value = ""
$ws.Range("AW12").Value = "Erland Lindblad"
$ws.Range("AX12").Value = "Jan Henriksson"
$ws.Range("AY12").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
